# Update "想去人数" (want-to-go count) values in column F
# for both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Row => new value for column F, shared by both sheets except row 13
$updates = @{
    3  = 2201
    4  = 87
    5  = 13107
    7  = 117
    8  = 516
    11 = 986
    12 = 13766
    22 = 1092
    25 = 5407
    26 = 939
    27 = 17
    28 = 314
    29 = 20
    30 = 40
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }

    # Row 13 has different source values per sheet but the same target value
    $ws.Cells.Item(13, 6).Value = 14357
}
